$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 42, shifting rows 42:57 down to 43:58
$ws.Rows(42).Insert()

# Populate the new row 42 with the new data point
$ws.Cells.Item(42, 1).Value = 11
$ws.Cells.Item(42, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42, 3).Value = "Bíobío"
$ws.Cells.Item(42, 4).Value = 45119
$ws.Cells.Item(42, 4).Style = $ws.Cells.Item(43, 4).Style
$ws.Cells.Item(42, 4).NumberFormat = $ws.Cells.Item(43, 4).NumberFormat
$ws.Cells.Item(42, 5).Value = 8
$ws.Cells.Item(42, 6).Value = 100112026
$ws.Cells.Item(42, 7).Value = "Haba"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 100
$ws.Cells.Item(42, 11).Value = 17000
$ws.Cells.Item(42, 12).Value = 18000
$ws.Cells.Item(42, 13).Value = 17500
$ws.Cells.Item(42, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(42, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(42, 16).Value = 700
$ws.Cells.Item(42, 17).Value = 25
$ws.Cells.Item(42, 18).Value = "Hortaliza"
